# SABolsas glossary updates
#  1. Merge the two "Avaliador" runs into a single run (no text change,
#     just a no-op Find/Replace on the exact text so the engine coalesces
#     the two adjacent, identically-formatted runs into one).
#  2. Insert a new blank paragraph + a new "SigPPG" glossary entry
#     paragraph right after the "Coordenador" paragraph.

$d = $word.ActiveDocument

# --- 1. Merge "Avaliador ... merece continuar ..." runs -------------------
$null = $d.Content.Find.Execute(
    "mesmo merece continuar", $true, $false, $false, $false, $false,
    $true, 1, $false, "mesmo merece continuar", 2)

# --- 2. Insert the two new paragraphs after the "Coordenador" paragraph ---
# Locate the paragraph that ends with "...podem acessar o sistema."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*podem acessar o sistema.*") {
        $target = $p
    }
}

$ns     = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rFonts = '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>'
$pPr    = "<w:pPr><w:pStyle w:val=""normal0""/><w:jc w:val=""both""/><w:rPr>$rFonts</w:rPr></w:pPr>"

# Empty paragraph
$newPara1 = "<w:p $ns>$pPr</w:p>"

# "SigPPG - Sistema Integrado de Pesquisa e Pos-graduacao." paragraph
$newPara2 = "<w:p $ns>$pPr" + `
    "<w:proofErr w:type=""spellStart""/>" + `
    "<w:r><w:rPr>$rFonts</w:rPr><w:t>SigPPG</w:t></w:r>" + `
    "<w:proofErr w:type=""spellEnd""/>" + `
    "<w:r><w:rPr>$rFonts</w:rPr><w:t xml:space=""preserve""> – Sistema Integrado</w:t></w:r>" + `
    "<w:r><w:rPr>$rFonts</w:rPr><w:t xml:space=""preserve""> de Pesquisa e Pós-graduação.</w:t></w:r>" + `
    "</w:p>"

$insertionPoint = $d.Range($target.Range.End, $target.Range.End)
$insertionPoint.InsertXML($newPara1 + $newPara2)

Write-Output "Edits applied"
